$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Update the existing "Layout" text (used by C8) to "Layout / Rode diesel pagina"
$ws.Range("C8").Value = "Layout / Rode diesel pagina"

# 2. Fill in the new row 9 data (date, hours, description)
$ws.Range("A9").Value = (Get-Date -Year 2015 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Layout / Rode diesel pagina / Website tijdelijk online zetten"

# 3. Extend the Totaal formula in D2 to include the new row
$ws.Range("D2").Formula = "=SUM(B2:B9)"

# 4. Update the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
